# Split the single BME_NSCG_LGH sheet into four site-specific sheets
# (BME_ARH, BME_BUH, BME_SMH, BME_CGH), each holding the same asset line
# (MONITORS, VIDEO) but keyed to a different Shop/Site so duplicate asset
# descriptions across sites no longer collide.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rename the original sheet to the first site.
$ws1.Name = "BME_ARH"

# Create the additional site sheets by copying the (now renamed) base
# sheet so formatting / conditional formatting / column widths carry over.
$ws1.Copy([System.Type]::Missing, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "BME_BUH"

$ws2.Copy([System.Type]::Missing, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "BME_SMH"

$ws3.Copy([System.Type]::Missing, $ws3)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "BME_CGH"

# ---- Sheet 1: BME_ARH ----
$ws1.Range("A1").Value = "BME_ARH: Annual Service Delivery Costs for Net New Equipment"
$ws1.Range("B3").Value = 47.58640699797898
$ws1.Range("B4").Value = 593735.6001137837
$ws1.Range("B5").Value = 307575.7425
$ws1.Range("B6").Value = 220540.37508
$ws1.Range("B7").Value = 65619.48253378378

$ws1.Range("A10").Value = "FHA"
$ws1.Range("B10").Value = "ARHCC"
$ws1.Range("C10").Value = "ARH"
$ws1.Range("D10").Value = "MONITORS, VIDEO"
$ws1.Range("E10").Value = 1
$ws1.Range("F10").Value = 3.542776494739543
$ws1.Range("G10").Value = 138.5225609443162
$ws1.Range("H10").Value = 307.1105651258654
$ws1.Range("I10").Value = 307.1105651258654

# ---- Sheet 2: BME_BUH ----
$ws2.Range("A1").Value = "BME_BUH: Annual Service Delivery Costs for Net New Equipment"
$ws2.Range("B3").Value = 54.16393726247094
$ws2.Range("B4").Value = 417502.3989112838
$ws2.Range("B5").Value = 209709.11
$ws2.Range("B6").Value = 142173.8063775
$ws2.Range("B7").Value = 65619.48253378378

$ws2.Range("A10").Value = "FHA"
$ws2.Range("B10").Value = "LMH_P"
$ws2.Range("C10").Value = "BHHO"
$ws2.Range("D10").Value = "MONITORS, VIDEO"
$ws2.Range("E10").Value = 1
$ws2.Range("F10").Value = 3.542776494739543
$ws2.Range("G10").Value = 142.8801760328458
$ws2.Range("H10").Value = 334.7708998288751
$ws2.Range("I10").Value = 334.7708998288751

# ---- Sheet 3: BME_SMH ----
$ws3.Range("A1").Value = "BME_SMH: Annual Service Delivery Costs for Net New Equipment"
$ws3.Range("B3").Value = 46.56042742184552
$ws3.Range("B4").Value = 1155583.248182784
$ws3.Range("B5").Value = 644836.015
$ws3.Range("B6").Value = 445127.7506489999
$ws3.Range("B7").Value = 65619.48253378378

$ws3.Range("A10").Value = "FHA"
$ws3.Range("B10").Value = "JPO"
$ws3.Range("C10").Value = "CAFVC"
$ws3.Range("D10").Value = "MONITORS, VIDEO"
$ws3.Range("E10").Value = 1
$ws3.Range("F10").Value = 3.542776494739543
$ws3.Range("G10").Value = 139.7935320118039
$ws3.Range("H10").Value = 304.7467198669448
$ws3.Range("I10").Value = 304.7467198669448

# ---- Sheet 4: BME_CGH ----
$ws4.Range("A1").Value = "BME_CGH: Annual Service Delivery Costs for Net New Equipment"
$ws4.Range("B3").Value = 146.324704283904
$ws4.Range("B4").Value = 228211.6669187838
$ws4.Range("B5").Value = 135024.6375
$ws4.Range("B6").Value = 27567.546885
$ws4.Range("B7").Value = 65619.48253378378

$ws4.Range("A10").Value = "FHA"
$ws4.Range("B10").Value = "RCH_C"
$ws4.Range("C10").Value = "CGH"
$ws4.Range("D10").Value = "MONITORS, VIDEO"
$ws4.Range("E10").Value = 1
$ws4.Range("F10").Value = 3.542776494739543
$ws4.Range("G10").Value = 138.5225609443162
$ws4.Range("H10").Value = 656.9182838810458
$ws4.Range("I10").Value = 656.9182838810458

$ws1.Activate()
